$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07785533333333333
$ws.Range("H2").Value = 0.233566
$ws.Range("I2").Value = 0.08237997085243232
$ws.Range("J2").Value = 0.08237997085243232
$ws.Range("M2").Value = 6.875726333333334
$ws.Range("N2").Value = 20.627179
$ws.Range("O2").Value = 0.6245871044820662
$ws.Range("P2").Value = 0.6245871044820662
$ws.Range("Q2").Value = 0.5353119655904445
$ws.Range("R2").Value = 4.817807690314
$ws.Range("S2").Value = 0.05145346746203772
$ws.Range("T2").Value = 0.05145346746203772

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07785533333333333
$ws.Range("H3").Value = 0.233566
$ws.Range("I3").Value = 0.08237997085243232
$ws.Range("J3").Value = 0.08237997085243232
$ws.Range("O3").Value = 0.06694469792011602
$ws.Range("P3").Value = 0.06694469792011602
$ws.Range("Q3").Value = 0.05737598098377777
$ws.Range("R3").Value = 0.5163838288539999
$ws.Range("S3").Value = 0.005514902263384044
$ws.Range("T3").Value = 0.005514902263384044

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07785533333333333
$ws.Range("H4").Value = 0.233566
$ws.Range("I4").Value = 0.08237997085243232
$ws.Range("J4").Value = 0.08237997085243232
$ws.Range("O4").Value = 0.3084681975978177
$ws.Range("P4").Value = 0.3084681975978177
$ws.Range("Q4").Value = 0.2643774038773333
$ws.Range("R4").Value = 2.379396634896001
$ws.Range("S4").Value = 0.02541160112701056
$ws.Range("T4").Value = 0.02541160112701056

$ws.Range("I5").Value = 0.8547824019796645
$ws.Range("J5").Value = 0.8547824019796644
$ws.Range("M5").Value = 6.875726333333334
$ws.Range("N5").Value = 20.627179
$ws.Range("O5").Value = 0.6245871044820662
$ws.Range("P5").Value = 0.6245871044820662
$ws.Range("Q5").Value = 5.554447798670779
$ws.Range("R5").Value = 49.99003018803701
$ws.Range("S5").Value = 0.5338860654147042
$ws.Range("T5").Value = 0.5338860654147042

$ws.Range("I6").Value = 0.8547824019796645
$ws.Range("J6").Value = 0.8547824019796644
$ws.Range("O6").Value = 0.06694469792011602
$ws.Range("P6").Value = 0.06694469792011602
$ws.Range("S6").Value = 0.05722314968795982
$ws.Range("T6").Value = 0.05722314968795982

$ws.Range("I7").Value = 0.8547824019796645
$ws.Range("J7").Value = 0.8547824019796644
$ws.Range("O7").Value = 0.3084681975978177
$ws.Range("P7").Value = 0.3084681975978177
$ws.Range("R7").Value = 24.68884547776801
$ws.Range("S7").Value = 0.2636731868770004
$ws.Range("T7").Value = 0.2636731868770004

$ws.Range("G8").Value = 0.05938633333333334
$ws.Range("I8").Value = 0.06283762716790325
$ws.Range("J8").Value = 0.06283762716790325
$ws.Range("M8").Value = 6.875726333333334
$ws.Range("N8").Value = 20.627179
$ws.Range("O8").Value = 0.6245871044820662
$ws.Range("P8").Value = 0.6245871044820662
$ws.Range("Q8").Value = 0.4083241759401112
$ws.Range("R8").Value = 3.674917583461001
$ws.Range("S8").Value = 0.03924757160532431
$ws.Range("T8").Value = 0.03924757160532431

$ws.Range("G9").Value = 0.05938633333333334
$ws.Range("I9").Value = 0.06283762716790325
$ws.Range("J9").Value = 0.06283762716790325
$ws.Range("O9").Value = 0.06694469792011602
$ws.Range("P9").Value = 0.06694469792011602
$ws.Range("Q9").Value = 0.04376513446344444
$ws.Range("R9").Value = 0.393886210171
$ws.Range("S9").Value = 0.004206645968772159
$ws.Range("T9").Value = 0.004206645968772159

$ws.Range("G10").Value = 0.05938633333333334
$ws.Range("I10").Value = 0.06283762716790325
$ws.Range("J10").Value = 0.06283762716790325
$ws.Range("O10").Value = 0.3084681975978177
$ws.Range("P10").Value = 0.3084681975978177
$ws.Range("Q10").Value = 0.2016612601893334
$ws.Range("S10").Value = 0.01938340959380678
$ws.Range("T10").Value = 0.01938340959380678

Write-Output "Applied TPM update to LR-pairs sheet"
